$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.398804783821106
$ws.Range("B1").Value = 3.847963809967041
$ws.Range("C1").Value = 3.389320850372314
$ws.Range("D1").Value = 3.649279117584229
$ws.Range("E1").Value = 1.300790548324585
